# Update cryptocurrency price / 1h-volume figures to the latest scrape.
# (GitHub Actions refresh job data.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.642.43'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '1.598.42'
$ws.Range("E3").Value = '  +0.12%  '

$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").Value = "'211.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("E6").Value = '  +0.44%  '

$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").Value = "'0.0618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("E9").Value = '  +0.13%  '

$ws.Range("D10").Value = "'19.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.41%  '

$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").Value = '1.822.99'
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").Value = '1.594.43'
$ws.Range("E13").Value = '  -0.37%  '

$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("E15").Value = '  +0.10%  '

$ws.Range("D16").Value = "'65.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.26%  '

$ws.Range("D17").Value = '26.646.64'
$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("D20").Value = "'208.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.72%  '

$ws.Range("D21").Value = "'7.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.48%  '

$ws.Range("D22").Value = "'4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("E23").Value = '  +0.95%  '

$ws.Range("D24").Value = "'8.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("D25").Value = "'145.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.97%  '

$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("E27").Value = '  -0.42%  '

$ws.Range("E28").Value = '  -0.52%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("E30").Value = '  +1.88%  '

$ws.Range("D31").Value = "'1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("E33").Value = '  +1.30%  '

$ws.Range("D34").Value = '1.277.86'
$ws.Range("E34").Value = '  -1.36%  '

$ws.Range("E35").Value = '  -9.47%  '

$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("E37").Value = '  +0.57%  '

$ws.Range("E38").Value = '  -0.66%  '

$ws.Range("E39").Value = '  -0.55%  '

$ws.Range("D40").Value = "'1.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +15.41%  '

$ws.Range("E41").Value = '  +2.88%  '

$ws.Range("D42").Value = "'2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.45%  '

$ws.Range("E43").Value = '  -0.78%  '

$ws.Range("D44").Value = "'63.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.45%  '

$ws.Range("D45").Value = '1.735.01'
$ws.Range("E45").Value = '  +0.20%  '

$ws.Range("D46").Value = "'90.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("E47").Value = '  -2.60%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = "'0.101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.86%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.0508"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.96%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.04%  '
